# Applies the OOXML diff: removes stray inner spaces / punctuation noise
# from several shared-string cell values across multiple worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "建物" (building / real estate) ---
$ws1 = $wb.Worksheets.Item("建物")
$ws1.Range("B2").Value = "臺北市大同區市府段一小段00927000建號"
$ws1.Range("D2").Value = "10000分之1"
$ws1.Range("F2").Value = "98年09月02日"
$ws1.Range("H2").Value = "5000000(為地上權房屋僅有使用權45年無土地所有權無權利持分）"

# --- Sheet "汽車" (car) ---
$ws2 = $wb.Worksheets.Item("汽車")
$ws2.Range("B2").Value = "TOYOTARAV4"
$ws2.Range("E2").Value = "102年02月03曰"

# --- Sheet "保險" (insurance) ---
$ws5 = $wb.Worksheets.Item("保險")
$ws5.Range("C3").Value = "富邦人壽真安心醫療養老保險"

# --- Sheet "債務" (debt) ---
$ws6 = $wb.Worksheets.Item("債務")
$ws6.Range("D2").Value = "京城銀行忠孝分行臺北市南港區忠孝東路六段21號"
$ws6.Range("F2").Value = "98年10月27日"
